# Updates meta_avg (H), meta_std (I) and meta_min (J) statistics
# for the "sub" histogram rows, per commit:
# "geração de histogramas para aon, flex e sub nas dimensões total, meta e contrib"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 1714.950834610086
$ws.Range("I2").Value = 2425.306729108642
$ws.Range("J2").Value = 0

# Row 3
$ws.Range("H3").Value = 937.0447030710235
$ws.Range("I3").Value = 1411.495693226725
$ws.Range("J3").Value = 0

# Row 4
$ws.Range("H4").Value = 2030.373249810451
$ws.Range("I4").Value = 4490.241197221475

# Row 5
$ws.Range("H5").Value = 1158.708690105903
$ws.Range("I5").Value = 2180.270411240455
$ws.Range("J5").Value = 0

# Row 6
$ws.Range("H6").Value = 696.4326100931312
$ws.Range("I6").Value = 764.1592346553673
$ws.Range("J6").Value = 0

# Row 7
$ws.Range("H7").Value = 552.7581271322479
$ws.Range("I7").Value = 791.154255909428
$ws.Range("J7").Value = 0

# Row 8
$ws.Range("H8").Value = 437.175996461933
$ws.Range("I8").Value = 1015.099894801972

# Row 9
$ws.Range("H9").Value = 638.8451893612514
$ws.Range("I9").Value = 1065.839823504247
$ws.Range("J9").Value = 0

# Row 10
$ws.Range("H10").Value = 780.2749801801036
$ws.Range("I10").Value = 1489.118149490151
$ws.Range("J10").Value = 0

# Row 13
$ws.Range("H13").Value = 756.444229088346
$ws.Range("I13").Value = 1173.798697147572
$ws.Range("J13").Value = 0

# Row 15
$ws.Range("H15").Value = 1182.830277432055
$ws.Range("I15").Value = 1570.66542202022
$ws.Range("J15").Value = 0

# Row 16
$ws.Range("H16").Value = 1186.777341475246
$ws.Range("I16").Value = 1442.796972601963

# Row 18
$ws.Range("H18").Value = 2300.173275326533
$ws.Range("I18").Value = 4921.295860260982
$ws.Range("J18").Value = 0

# Row 20
$ws.Range("H20").Value = 1129.303003825498
$ws.Range("I20").Value = 1371.580839412833

# Row 21
$ws.Range("H21").Value = 932.2401095754986
$ws.Range("I21").Value = 1121.464969602491
$ws.Range("J21").Value = 0

# Row 22
$ws.Range("H22").Value = 297.8228084596681
$ws.Range("I22").Value = 369.1913163392061
